$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the shared string "He-Regular" -> "He-Regular1" everywhere it is used
# (cells C4 and C5 both reference it).
$ws.Range("C4").Value = "He-Regular1"
$ws.Range("C5").Value = "He-Regular1"

# Update the saved cursor/selection position.
$ws.Range("C15").Select()
